$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

$ws.Range("D2").Value = 3809.51
$ws.Range("E2").Value = -3809.51

$ws.Range("D4").Value = 4131.360000000001
$ws.Range("E4").Value = 13368.64
$ws.Range("F4").Value = 0.2360777142857143
